$wb = $excel.ActiveWorkbook

# Sheet 1: "Datos del Cliente" -> row 2 (A2:H2) all become "N/A"
$ws1 = $wb.Worksheets.Item("Datos del Cliente")
$ws1.Range("A2:H2").Value = "N/A"

# Sheet 2: "Datos de Lectura" -> row 2 (A2:E2) all become "N/A"
$ws2 = $wb.Worksheets.Item("Datos de Lectura")
$ws2.Range("A2:E2").Value = "N/A"

# Sheet 3: "Costos de Energía" -> row 2 (A2:H2) all become "N/A"
$ws3 = $wb.Worksheets.Item("Costos de Energía")
$ws3.Range("A2:H2").Value = "N/A"

# Sheet 4: "Desglose de Importe" -> row 2 (A2:G2) all become "N/A"
$ws4 = $wb.Worksheets.Item("Desglose de Importe")
$ws4.Range("A2:G2").Value = "N/A"

# Sheet 5: "Consumo Histórico" -> clear entire sheet content (headers + all data rows)
$ws5 = $wb.Worksheets.Item("Consumo Histórico")
$ws5.Cells.Clear()
